# "aggiornamento fino a 6/03" - append 3 new daily rows (245-247) to Sheet1,
# continuing the existing date series (serial dates 44319, 44320, 44321)
# with the same values as the preceding row (244): nuovi pos. = 0,
# somma mobile 7gg. = 1, somma mobile 7gg. per 100mila abitanti = 28.87669650591972.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 244
$newDates = @(44319, 44320, 44321)
$newB = 0
$newC = 1
$newD = 28.87669650591972

$row = $lastRow + 1
foreach ($d in $newDates) {
    # Reuse the exact formatting (date number format, bold font, border,
    # alignment) already applied to column A of the last existing row.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = $newB
    $ws.Cells.Item($row, 3).Value = $newC
    $ws.Cells.Item($row, 4).Value = $newD

    $row = $row + 1
}
